$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.087.22'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.780.24'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.46'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.78'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.55%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.50%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.036.71'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.783.83'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.062.30'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.62'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.65'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0788'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.83%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.32%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.68%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.56'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.63%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.60%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.46%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.26%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.446.59'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.54%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +4.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.650'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.13%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.40'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.72%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.914'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.68'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0519'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.52%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.08'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.24%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.04'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.937.75'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '104.17'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.23%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.63%  '
